$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Test Case Name for row 4 (Framework_002 -> Framework_003)
$ws.Range("A4").Value = "Framework_003"

# Duplicate row 4 into row 5, carrying over values + formatting/styles exactly
$ws.Range("A4:M4").Copy($ws.Range("A5:M5"))

# Overwrite the values that differ between row 4 and the new row 5
$ws.Range("A5").Value = "Framework_003"
$ws.Range("B5").Value = "testuser_1"
$ws.Range("C5").Value = "Test@123"
$ws.Range("D5").Value = "Chrome"
$ws.Range("E5").Value = "iMacs"
$ws.Range("F5").Value = "Product 1"
$ws.Range("G5").Value = "Pratik"
$ws.Range("H5").Value = "Sharma"
$ws.Range("I5").Value = "29, LimeSquare, City Road"
$ws.Range("J5").Value = "Newcastle"
$ws.Range("K5").Value = "United Kingdom"
$ws.Range("L5").Value = "07777777777"
$ws.Range("M5").Value = "tooolsqa@gmail.com"

# Recreate the hyperlinks on the new row (mirrors C/M hyperlinks used elsewhere)
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("M5"), "mailto:tooolsqa@gmail.com")

# Extend the data validation ranges to include the new row
$ws.Range("E2:E5").Validation.Delete()
$ws.Range("E2:E5").Validation.Add(3, 1, 1, "Accessories, iMacs, iPads, iPhones")
$ws.Range("E2:E5").Validation.InputMessage = ""
$ws.Range("E2:E5").Validation.ShowInput = $true
$ws.Range("E2:E5").Validation.ShowError = $true

$ws.Range("F2:F5").Validation.Delete()
$ws.Range("F2:F5").Validation.Add(3, 1, 1, "Product 1, Product 2, Product 3, Product 4")
$ws.Range("F2:F5").Validation.InputMessage = ""
$ws.Range("F2:F5").Validation.ShowInput = $true
$ws.Range("F2:F5").Validation.ShowError = $true

# Update the view: scroll so column B is the left-most visible column, and select I9
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("I9").Select()
